$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spike trap (尖刺) effect text: changed from reducing rank by 1, to reducing
# rank by 1 and also sending the card to graveyard if the row ends up empty.
$ws.Range("D5").Value = "交锋时：同一行中所有怪物牌点数点数减1。该效果结算完毕后，如果本行没有怪物牌，则将本牌也送墓。"

# Update the sheet's active selection to D6, matching the saved workbook state.
$ws.Range("D6").Select()
